$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    15  = @("ba", "Appreciation")
    40  = @("sd", "Statement-non-opinion")
    64  = @("sv", "Statement-opinion")
    85  = @("sv", "Statement-opinion")
    86  = @("b", "Acknowledge (Backchannel)")
    99  = @("sv", "Statement-opinion")
    101 = @("sv", "Statement-opinion")
    103 = @("b", "Acknowledge (Backchannel)")
    106 = @("b", "Acknowledge (Backchannel)")
    121 = @("aa", "Agree/Accept")
    122 = @("sd", "Statement-non-opinion")
    123 = @("sd", "Statement-non-opinion")
    125 = @("sv", "Statement-opinion")
    126 = @("sv", "Statement-opinion")
    133 = @("b", "Acknowledge (Backchannel)")
    140 = @("sd", "Statement-non-opinion")
    149 = @("aa", "Agree/Accept")
    151 = @("sd", "Statement-non-opinion")
    159 = @("b", "Acknowledge (Backchannel)")
    167 = @("ba", "Appreciation")
    173 = @("aa", "Agree/Accept")
    176 = @("sd", "Statement-non-opinion")
    181 = @("sd", "Statement-non-opinion")
    190 = @("aa", "Agree/Accept")
    193 = @("aa", "Agree/Accept")
    194 = @("b", "Acknowledge (Backchannel)")
    197 = @("b", "Acknowledge (Backchannel)")
    200 = @("aa", "Agree/Accept")
    201 = @("aa", "Agree/Accept")
    204 = @("b", "Acknowledge (Backchannel)")
    226 = @("sd", "Statement-non-opinion")
    300 = @("sd", "Statement-non-opinion")
    305 = @("sd", "Statement-non-opinion")
    312 = @("sd", "Statement-non-opinion")
    314 = @("aa", "Agree/Accept")
    317 = @("sd", "Statement-non-opinion")
    329 = @("sd", "Statement-non-opinion")
    336 = @("sd", "Statement-non-opinion")
    342 = @("sv", "Statement-opinion")
    350 = @("sd", "Statement-non-opinion")
    351 = @("ba", "Appreciation")
    358 = @("sd", "Statement-non-opinion")
    365 = @("ba", "Appreciation")
    367 = @("sv", "Statement-opinion")
    383 = @("aa", "Agree/Accept")
    391 = @("b", "Acknowledge (Backchannel)")
    393 = @("sd", "Statement-non-opinion")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
